# Append daily records (9 Dec 2021 - 6 Jan 2022) to the Frassinoro report sheet.
# Matches commit: "aggiornamento fino a 6 gennaio 2022"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows: date-serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$data = @(
    @(44539,0,0,0),
    @(44540,0,0,0),
    @(44541,0,0,0),
    @(44542,0,0,0),
    @(44543,0,0,0),
    @(44544,1,1,53.73455131649651),
    @(44545,0,1,53.73455131649651),
    @(44546,0,1,53.73455131649651),
    @(44547,0,1,53.73455131649651),
    @(44548,2,3,161.2036539494895),
    @(44550,1,4,214.9382052659861),
    @(44551,0,4,214.9382052659861),
    @(44552,0,3,161.2036539494895),
    @(44553,0,3,161.2036539494895),
    @(44554,1,4,214.9382052659861),
    @(44555,1,5,268.6727565824825),
    @(44556,3,6,322.4073078989791),
    @(44557,4,9,483.6109618484685),
    @(44558,3,12,644.8146157979581),
    @(44559,4,16,859.7528210639442),
    @(44560,5,21,1128.425577646427),
    @(44561,7,27,1450.832885545406),
    @(44562,5,31,1665.771090811392),
    @(44563,1,29,1558.301988178399),
    @(44564,0,25,1343.363782912413),
    @(44565,1,23,1235.89468027942),
    @(44566,1,20,1074.69102632993)
)

$lastExistingRow = 464
$startRow = $lastExistingRow + 1
$rowCount = $data.Count
$colCount = 4
$endRow = $startRow + $rowCount - 1

$targetRange = "A$startRow`:D$endRow"

# Copy formatting (style, number format, borders, etc.) from the last existing
# data row down across every new row before writing the values.
$ws.Range("A$lastExistingRow`:D$lastExistingRow").Copy()
$ws.Range($targetRange).PasteSpecial(-4122)  # xlPasteFormats

# Build a 2D array to write all new cell values in one shot.
$arr = New-Object 'object[,]' $rowCount,$colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($j = 0; $j -lt $colCount; $j++) {
        $arr[$i,$j] = $data[$i][$j]
    }
}

$ws.Range($targetRange).Value = $arr
